# Applies the "reformat architecture and update the tester nodes" edit to
# report_TensorrtExecutionProvider.xlsx
#
# Summary of changes:
#  - Rows 9, 54, 123 on "Report": used_provider (C) gets filled in with
#    "TensorrtExecutionProvider" and status (D) flips from a FAIL message to
#    "SUCCESS (via decomposition)" (green fill, matching other such cells).
#  - Rows 15, 29, 55, 56, 88, 89, 127, 137 on "Report": status (D) text
#    changes from "SUCCESS WITH FALLBACK" to "SUCCESS (via decomposition)".
#  - "Data_PieChart" summary sheet: SUCCESS / FAIL roll-up counts and
#    percentages are refreshed (+3 SUCCESS, -3 FAIL) and the report
#    timestamp is updated.

$wb = $excel.ActiveWorkbook

$report = $wb.Worksheets.Item("Report")
$pie    = $wb.Worksheets.Item("Data_PieChart")

# Green fill color used by existing "SUCCESS*" cells (RGB 00AA44 -> BGR long)
$successGreen = 4499968

# --- Rows that moved from a hard FAIL to SUCCESS (via decomposition) -------
$newSuccessRows = @(9, 54, 123)
foreach ($r in $newSuccessRows) {
    $report.Range("C$r").Value = "TensorrtExecutionProvider"
    $report.Range("D$r").Value = "SUCCESS (via decomposition)"
    $report.Range("D$r").Interior.Color = $successGreen
}

# --- Rows whose status text changed from WITH FALLBACK to via decomposition
$relabelRows = @(15, 29, 55, 56, 88, 89, 127, 137)
foreach ($r in $relabelRows) {
    $report.Range("D$r").Value = "SUCCESS (via decomposition)"
}

# --- Refresh the summary pie-chart data sheet ------------------------------
$pie.Range("B2").Value = 232
$pie.Range("C2").Value = 88.5

$pie.Range("B9").Value = 21
$pie.Range("C9").Value = 8

$pie.Range("B10").Value = "2025-11-18 14:41:55"
